# adding averages and more checks
$wb = $excel.ActiveWorkbook

$training = $wb.Worksheets.Item("Training Dashboard")
$exam = $wb.Worksheets.Item("Exam Dashboard")

# --- Training Dashboard: update PERIOD TO EXPIRE (col H) and LAST UPDATE (col I) for rows 3-32 ---
$periods = @{
    3 = 676; 4 = 368; 5 = 338; 6 = 334; 7 = 357; 8 = 329; 9 = 355; 10 = 370;
    11 = 697; 12 = 697; 13 = 338; 14 = 377; 15 = 490; 16 = 489; 17 = 489; 18 = 489;
    19 = 47; 20 = 298; 21 = 298; 22 = -42; 23 = 300; 24 = 204; 25 = 205; 26 = 311;
    27 = 304; 28 = 311; 29 = 312; 30 = 311; 31 = 332; 32 = 332
}

foreach ($row in 3..32) {
    $training.Cells.Item($row, 8).Value = $periods[$row]
    $training.Cells.Item($row, 9).Value = "16-Sep-2025"
}

# --- Exam Dashboard: update COMMENTS (col E) for rows 3-7 ---
foreach ($row in 3..7) {
    $exam.Cells.Item($row, 5).Value = "date is valid"
}

# narrow column E (was widened for old long comments)
$exam.Columns.Item(5).ColumnWidth = 15

# --- styles.xml: header row font becomes bold white (on dark blue fill); drop the 14pt title font size ---
$headerRange1 = $training.Range("A2:K2")
$headerRange1.Font.Bold = $true
$headerRange1.Font.Color = "#FFFFFF"

$headerRange2 = $exam.Range("A2:G2")
$headerRange2.Font.Bold = $true
$headerRange2.Font.Color = "#FFFFFF"

$training.Range("A1").Font.Size = 11
$training.Range("A1").Font.Color = "#FFFFFF"

$exam.Range("A1").Font.Size = 11
$exam.Range("A1").Font.Color = "#FFFFFF"
